$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything currently in A:W shifts to B:X.
$ws.Columns("A:A").Insert()

# New column A = "Match ID". Header + data rows get the bold "no border" style
# (same font as the other header style, just without the border/alignment).
$ws.Range("A1:A19").Font.Bold = $true

$ws.Range("A1").Value = "Match ID"

# Data rows (4-19) all belong to the same match.
$ws.Range("A4:A19").Value = 17

# Row 20 is a hidden totals row. Writing into a hidden row makes the engine
# stamp an explicit row height, so temporarily unhide it, write the value,
# then hide it again to match the original (default-height) row.
$ws.Rows("20:20").Hidden = $false
$ws.Range("A20").Value = 17
$ws.Rows("20:20").Hidden = $true

# Update the selection to mirror the authored file (A1:A19 selected).
$ws.Range("A1:A19").Select() | Out-Null
